$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 54; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E = purpose
    if ($cell.Value() -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
